$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new column F "SourceCodeManagement" with a couple of sample rows,
# matching the new checkbox method used for week 3 data-driven testing.
$ws.Range("F1").Value = "SourceCodeManagement"
$ws.Range("F3").Value = "CVS"
$ws.Range("F2").Value = "None"

# Copy header/body formatting from column E onto the new column F cells
# so the new column matches the look (bold header, text format) of the
# rest of the table.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("E2:E3").Copy()
$ws.Range("F2:F3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match column F's width to column E's width.
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(5).ColumnWidth()

# Update the active selection on the sheet.
$null = $ws.Range("C5").Select()
